# "access data for report"
# - rename "enrolments" sheet to "invoices" (cascades: defined-name filter ref)
# - fix header casing: customer_Id -> customer_id, instructor_Id -> instructor_id,
#   schedule_Id -> schedule_id, enrolment_Id -> invoice_id
# - schedule sheet: class_Id header becomes class_id (matching the class_id label
#   already used on the classes sheet)
# - selection / active-sheet bookkeeping left by the editor moving around the workbook

$wb = $excel.ActiveWorkbook

# Rename the "enrolments" sheet to "invoices" (also updates the
# _xlnm._FilterDatabase defined name that points at it).
$invoices = $wb.Worksheets.Item("enrolments")
$invoices.Name = "invoices"

# Normalize header labels across the sheets that reference these columns.
$customers = $wb.Worksheets.Item("customers")
$customers.Range("A1").Value = "customer_id"

$invoices.Range("A1").Value = "invoice_id"
$invoices.Range("B1").Value = "customer_id"
$invoices.Range("C1").Value = "schedule_id"

$instructors = $wb.Worksheets.Item("instructors")
$instructors.Range("A1").Value = "instructor_id"

$schedule = $wb.Worksheets.Item("schedule")
$schedule.Range("A1").Value = "schedule_id"
$schedule.Range("B1").Value = "class_id"
$schedule.Range("C1").Value = "instructor_id"

# Leave the workbook with "customers" the active tab/selected sheet, and
# restore the per-sheet selections to where the editor left them.
$classes = $wb.Worksheets.Item("classes")
$classes.Range("C1").Select() | Out-Null

$invoices.Range("A1").Select() | Out-Null

$instructors.Range("E1").Select() | Out-Null

$schedule.Range("A1").Select() | Out-Null

$customers.Activate() | Out-Null
$customers.Range("E1").Select() | Out-Null
